# Scheduled market-data refresh for the Brynhildr Profits workbook.
#
# Each job sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) has a
# Table_<JOB> table spanning A1:N141 with per-leve crafting-profit
# math in columns H:N:
#   H = currentAveragePrice        K = LevePriceNQ
#   I = currentAveragePriceNQ      L = LevePriceHQ
#   J = currentAveragePriceHQ      M = LeveProfitNQ
#   N = LeveProfitHQ
#
# This runner only refreshes those derived market-price/profit
# cells (H:N) for the rows the latest market snapshot touched;
# leve metadata in A:G is untouched.

$wb = $excel.ActiveWorkbook

# Map of sheet name -> list of { Cell, Value } updates to apply.
$sheetUpdates = @{}

$sheetUpdates["ALC"] = @(
    @{ Cell = "H12"; Value = 147.84616 },
    @{ Cell = "I12"; Value = 141.1 },
    @{ Cell = "K12"; Value = 141.1 },
    @{ Cell = "M12"; Value = 28.90000000000001 },
    @{ Cell = "H86"; Value = 26449.5 },
    @{ Cell = "I86"; Value = 2899 },
    @{ Cell = "J86"; Value = 50000 },
    @{ Cell = "K86"; Value = 2899 },
    @{ Cell = "L86"; Value = 50000 },
    @{ Cell = "M86"; Value = -1776 },
    @{ Cell = "N86"; Value = -52246 },
    @{ Cell = "H89"; Value = 26449.5 },
    @{ Cell = "I89"; Value = 2899 },
    @{ Cell = "J89"; Value = 50000 },
    @{ Cell = "K89"; Value = 14495 },
    @{ Cell = "L89"; Value = 250000 },
    @{ Cell = "M89"; Value = -8879 },
    @{ Cell = "N89"; Value = -261232 },
    @{ Cell = "H101"; Value = 1727 },
    @{ Cell = "I101"; Value = 268 },
    @{ Cell = "J101"; Value = 3186 },
    @{ Cell = "K101"; Value = 804 },
    @{ Cell = "L101"; Value = 9558 },
    @{ Cell = "M101"; Value = 818 },
    @{ Cell = "N101"; Value = -12802 },
    @{ Cell = "H106"; Value = 12272.286 },
    @{ Cell = "J106"; Value = 14472 },
    @{ Cell = "L106"; Value = 14472 },
    @{ Cell = "N106"; Value = -15734 },
    @{ Cell = "H138"; Value = 2597.1 },
    @{ Cell = "I138"; Value = 3000.4614 },
    @{ Cell = "J138"; Value = 2160.125 },
    @{ Cell = "K138"; Value = 9001.3842 },
    @{ Cell = "L138"; Value = 6480.375 },
    @{ Cell = "M138"; Value = -3861.3842 },
    @{ Cell = "N138"; Value = -16760.375 }
)

$sheetUpdates["ARM"] = @(
    @{ Cell = "H5"; Value = 3184.5715 },
    @{ Cell = "J5"; Value = 4100.4 },
    @{ Cell = "L5"; Value = 4100.4 },
    @{ Cell = "N5"; Value = -4324.4 },
    @{ Cell = "H32"; Value = 206276.4 },
    @{ Cell = "I32"; Value = 231019.7 },
    @{ Cell = "K32"; Value = 231019.7 },
    @{ Cell = "M32"; Value = -230732.7 },
    @{ Cell = "H33"; Value = 37676.332 },
    @{ Cell = "I33"; Value = 36500 },
    @{ Cell = "K33"; Value = 36500 },
    @{ Cell = "M33"; Value = -36171 },
    @{ Cell = "H122"; Value = 1503.591 },
    @{ Cell = "I122"; Value = 1431.45 },
    @{ Cell = "K122"; Value = 4294.35 },
    @{ Cell = "M122"; Value = -1844.35 },
    @{ Cell = "H132"; Value = 4826.3193 },
    @{ Cell = "I132"; Value = 3289.3618 },
    @{ Cell = "K132"; Value = 9868.0854 },
    @{ Cell = "M132"; Value = -7338.0854 }
)

$sheetUpdates["BSM"] = @(
    @{ Cell = "H4"; Value = 3184.5715 },
    @{ Cell = "J4"; Value = 4100.4 },
    @{ Cell = "L4"; Value = 4100.4 },
    @{ Cell = "N4"; Value = -4330.4 }
)

$sheetUpdates["CRP"] = @(
    @{ Cell = "H6"; Value = 198125 },
    @{ Cell = "I6"; Value = 215000 },
    @{ Cell = "K6"; Value = 215000 },
    @{ Cell = "M6"; Value = -214887 },
    @{ Cell = "H7"; Value = 89.545456 },
    @{ Cell = "I7"; Value = 38.666668 },
    @{ Cell = "K7"; Value = 38.666668 },
    @{ Cell = "M7"; Value = 74.333332 },
    @{ Cell = "H58"; Value = 3127.4546 },
    @{ Cell = "I58"; Value = 3210.5 },
    @{ Cell = "J58"; Value = 3027.8 },
    @{ Cell = "K58"; Value = 3210.5 },
    @{ Cell = "L58"; Value = 3027.8 },
    @{ Cell = "M58"; Value = -3007.5 },
    @{ Cell = "N58"; Value = -3433.8 },
    @{ Cell = "H136"; Value = 3127.4546 },
    @{ Cell = "I136"; Value = 3210.5 },
    @{ Cell = "J136"; Value = 3027.8 },
    @{ Cell = "K136"; Value = 9631.5 },
    @{ Cell = "L136"; Value = 9083.400000000001 },
    @{ Cell = "M136"; Value = -7081.5 },
    @{ Cell = "N136"; Value = -14183.4 }
)

$sheetUpdates["CUL"] = @(
    @{ Cell = "H31"; Value = 4433.6665 },
    @{ Cell = "I31"; Value = 5150.5 },
    @{ Cell = "K31"; Value = 15451.5 },
    @{ Cell = "M31"; Value = -15163.5 },
    @{ Cell = "H50"; Value = 147518 },
    @{ Cell = "I50"; Value = 418.57895 },
    @{ Cell = "J50"; Value = 333843.94 },
    @{ Cell = "K50"; Value = 1255.73685 },
    @{ Cell = "L50"; Value = 1001531.82 },
    @{ Cell = "M50"; Value = -774.73685 },
    @{ Cell = "N50"; Value = -1002493.82 },
    @{ Cell = "H53"; Value = 147518 },
    @{ Cell = "I53"; Value = 418.57895 },
    @{ Cell = "J53"; Value = 333843.94 },
    @{ Cell = "K53"; Value = 1255.73685 },
    @{ Cell = "L53"; Value = 1001531.82 },
    @{ Cell = "M53"; Value = -774.73685 },
    @{ Cell = "N53"; Value = -1002493.82 },
    @{ Cell = "H121"; Value = 16647.176 },
    @{ Cell = "I121"; Value = 245.57143 },
    @{ Cell = "J121"; Value = 28128.3 },
    @{ Cell = "K121"; Value = 736.71429 },
    @{ Cell = "L121"; Value = 84384.89999999999 },
    @{ Cell = "M121"; Value = 573.28571 },
    @{ Cell = "N121"; Value = -87004.89999999999 },
    @{ Cell = "H129"; Value = 2272.7273 },
    @{ Cell = "J129"; Value = 2389 },
    @{ Cell = "L129"; Value = 7167 },
    @{ Cell = "N129"; Value = -17167 },
    @{ Cell = "H131"; Value = 2417.3462 },
    @{ Cell = "I131"; Value = 1149.75 },
    @{ Cell = "J131"; Value = 2522.9792 },
    @{ Cell = "K131"; Value = 3449.25 },
    @{ Cell = "L131"; Value = 7568.937600000001 },
    @{ Cell = "M131"; Value = 1590.75 },
    @{ Cell = "N131"; Value = -17648.9376 }
)

$sheetUpdates["GSM"] = @(
    @{ Cell = "H70"; Value = 22007.592 },
    @{ Cell = "I70"; Value = 20669.945 },
    @{ Cell = "K70"; Value = 20669.945 },
    @{ Cell = "M70"; Value = -20399.945 },
    @{ Cell = "H73"; Value = 22007.592 },
    @{ Cell = "I73"; Value = 20669.945 },
    @{ Cell = "K73"; Value = 20669.945 },
    @{ Cell = "M73"; Value = -19733.945 },
    @{ Cell = "H132"; Value = 13150.5 },
    @{ Cell = "I132"; Value = 26395 },
    @{ Cell = "J132"; Value = 4320.8335 },
    @{ Cell = "K132"; Value = 79185 },
    @{ Cell = "L132"; Value = 12962.5005 },
    @{ Cell = "M132"; Value = -76655 },
    @{ Cell = "N132"; Value = -18022.5005 }
)

$sheetUpdates["LTW"] = @(
    @{ Cell = "H68"; Value = 8076.5 },
    @{ Cell = "I68"; Value = 9649.916999999999 },
    @{ Cell = "J68"; Value = 4929.6665 },
    @{ Cell = "K68"; Value = 9649.916999999999 },
    @{ Cell = "L68"; Value = 4929.6665 },
    @{ Cell = "M68"; Value = -8900.916999999999 },
    @{ Cell = "N68"; Value = -6427.6665 },
    @{ Cell = "H71"; Value = 8076.5 },
    @{ Cell = "I71"; Value = 9649.916999999999 },
    @{ Cell = "J71"; Value = 4929.6665 },
    @{ Cell = "K71"; Value = 48249.585 },
    @{ Cell = "L71"; Value = 24648.3325 },
    @{ Cell = "M71"; Value = -44505.585 },
    @{ Cell = "N71"; Value = -32136.3325 }
)

$sheetUpdates["WVR"] = @(
    @{ Cell = "H122"; Value = 79553.53 },
    @{ Cell = "I122"; Value = 4671.091 },
    @{ Cell = "K122"; Value = 14013.273 },
    @{ Cell = "M122"; Value = -11563.273 },
    @{ Cell = "H126"; Value = 2174.7368 },
    @{ Cell = "I126"; Value = 2004.7333 },
    @{ Cell = "K126"; Value = 6014.199900000001 },
    @{ Cell = "M126"; Value = -3544.199900000001 },
    @{ Cell = "H136"; Value = 1695.15 },
    @{ Cell = "I136"; Value = 1392.7693 },
    @{ Cell = "K136"; Value = 4178.3079 },
    @{ Cell = "M136"; Value = -1628.3079 }
)

$totalUpdates = 0
foreach ($sheetName in $sheetUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($u in $sheetUpdates[$sheetName]) {
        $ws.Range($u.Cell).Value = $u.Value
        $totalUpdates = $totalUpdates + 1
    }
}

Write-Output "Updated $totalUpdates cells across $($sheetUpdates.Keys.Count) sheets."
